$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Journals sheet: publication rows 3-8 got reshuffled / updated as
# several "In Press" articles were officially published (new
# volume/issue/page/DOI info), and one brand-new entry was folded in.
# ------------------------------------------------------------------
$wsJ = $wb.Worksheets.Item("Journals")

# Row 3 (was Karlin "Impact of Real-World..." / In Press) ->
#        Stipancic et al, Autoscore (still In Press)
$wsJ.Range("A3").Value = 'Stipancic, K. L., Barrett, T. S., Tjaden, K., & Borrie, S. A.'
$wsJ.Range("C3").Value = 'Automated scoring of the Speech Intelligibility Test using Autoscore.'
$wsJ.Range("D3").Value = 'American Journal of Speech-Language Pathology'
$wsJ.Range("H3").Clear()
$wsJ.Rows.Item(3).RowHeight = 34

# Row 4 (was Stipancic Autoscore) -> Borrie et al, dysarthria training
$wsJ.Range("A4").Value = 'Borrie, S. A., Tetzloff, K., Barrett, T. S., & Lansford, K. L.'
$wsJ.Range("C4").Value = 'Increasing motivation increases intelligibility benefits of perceptual training of dysarthria.'

# Row 5 (was Popkov et al depression cost) -> Blais et al, stigma/MST,
# now has a DOI (online first, no volume/issue/page yet)
$wsJ.Range("A5").Value = 'Blais, R. K., Barrett, T. S., Tannahill, H. S., & Hoyt, T. '
$wsJ.Range("C5").Value = 'Fears of Stigma Perceived from Unit Leaders for Seeking Psychological Services for Military Sexual Trauma Sequelae is Associated with Suicide Risk among Male Service Members and Veterans.'
$wsJ.Range("D5").Value = 'Stigma and Health'
$wsJ.Range("H5").Value = 'doi: 10.1037/sah0000543'

# Row 6 (was Borrie et al dysarthria training) -> Popkov et al
# depression cost study, now officially published in 2025 with
# volume/page/DOI filled in
$wsJ.Range("A6").Value = 'Popkov, A. A., Barrett, T. S., Shergil, A., Donohue, M., Anderson, R. J., & Karlin, B. E.'
$wsJ.Range("B6").Value = 2025
$wsJ.Range("C6").Value = 'Association Between Depression Symptom Severity and Total Cost of Care: Findings from a Large, 2-year, Claims-Based, Retrospective Population Health Study.'
$wsJ.Range("D6").Value = 'Journal of Affective Disorders'
$wsJ.Range("E6").Value = 368
$wsJ.Range("G6").Value = '41-47'
$wsJ.Range("H6").Value = 'doi: 10.1016/j.jad.2024.09.056'
$wsJ.Rows.Item(6).RowHeight = 68

# Row 7 (was Blais et al stigma/MST) -> Karlin et al "Impact of
# Real-World Implementation..." now officially published with
# volume/issue/page/DOI (matches the zpae053 formatting of row 8)
$wsJ.Range("A7").Value = 'Karlin, B. E., Anderson, R. J., Rung, J. M., Drury-Gworek, C., & Barrett, T. S.'
$wsJ.Range("B7").Value = 2024
$wsJ.Range("C7").Value = 'Impact of Real-World Implementation of Evidence-Based Insomnia Treatment within a Large Payor-Provider Health System: Initial Provider and Patient-Level Outcomes.'
$wsJ.Range("D7").Value = 'SLEEP Advances'
$wsJ.Range("E7").Value = 5
$wsJ.Range("F7").Value = 1
$wsJ.Range("G8").Copy()
$wsJ.Range("G7").PasteSpecial(-4122)
$wsJ.Range("G7").Value = 'zpae053'
$wsJ.Range("H7").Value = 'doi: 10.1093/sleepadvances/zpae053'

# Row 8 (Clinical and Financial Significance...) keeps its content,
# just gains volume/issue info (5/1)
$wsJ.Range("E8").Value = 5
$wsJ.Range("F8").Value = 1

$wsJ.Activate()
$wsJ.Range("G5").Select()

# NOTE: Grants automatically loses its tabSelected flag once Journals
# becomes the active sheet above, and its own selection (E3) is left
# untouched since we never re-select anything on that sheet.
